$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new log row ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A6").Value = "Opvolging klacht"
$logs.Range("B6").Value = "mailmind.test@zohomail.eu"
$logs.Range("D6").Value = "Klacht / Probleem"
$logs.Range("F6").Value = "2025-08-30 19:05:06"
$logs.Range("G6").Value = "Nee"
$logs.Range("H6").Value = "Ja"
$logs.Range("I6").Value = "Nee"
$logs.Range("J6").Value = "Nee"

# Extend the conditional-formatting ranges to include the new row 6
$colRanges = @("D2:D5", "G2:G5", "H2:H5", "I2:I5", "J2:J5")
foreach ($oldRange in $colRanges) {
    $col = $oldRange.Substring(0, 1)
    $newRange = "$col" + "2:$col" + "6"
    $fcs = $logs.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# --- Dashboard sheet: refresh the category counts (Klacht / Probleem now 2) ---
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A3").Value = "Klacht / Probleem"
$dashboard.Range("B3").Value = 2
$dashboard.Range("A4").Value = "Planning / Afspraak"
$dashboard.Range("B4").Value = 1
